# Remove column M from the alcohol data sheet (Sheet1), shifting column N left.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Columns("M").Delete()

$ws.Range("M1").Select()
